$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview": add a row for the new c131e10c-... handoff file.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$rowOverview = $loOverview.ListRows.Add()

$wsOverview.Range("A10").Value = "c131e10c-ac85-432b-a277-ff971a1a9a8b.md"
$wsOverview.Hyperlinks.Add($wsOverview.Range("B10"), "https://github.com/OpenLocalizationTestOrg/ol-test4/blob/c131e10c1e1c131e10c1e1c131e10c1e1c131e10/e2e/c131e10c-ac85-432b-a277-ff971a1a9a8b.md", "", "", "e2e\c131e10c-ac85-432b-a277-ff971a1a9a8b.md")
$wsOverview.Range("C10").Value = ".md"
$wsOverview.Range("D10").Value = "'"
$wsOverview.Range("D10").Style = "Normal"
$wsOverview.Range("E10").Value = "Ready for handoff"
$wsOverview.Range("F10").Value = "Ready for handoff"
$wsOverview.Range("G10").Value = "2017-02-21 02:59:17"
$wsOverview.Range("G10").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------------
# Sheet "zh-cn": add the matching detail row.
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)
$rowZhCn = $loZhCn.ListRows.Add()

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A10"), "https://github.com/OpenLocalizationTestOrg/ol-test4/blob/c131e10c1e1c131e10c1e1c131e10c1e1c131e10/e2e/c131e10c-ac85-432b-a277-ff971a1a9a8b.md", "", "", "c131e10c-ac85-432b-a277-ff971a1a9a8b.md")
$wsZhCn.Range("B10").Value = ".md"
$wsZhCn.Range("C10").Value = "Ready for handoff"
$wsZhCn.Range("D10").Value = "e2e"
$wsZhCn.Range("E10").Value = "ht"
$wsZhCn.Range("F10").Value = "False"
$wsZhCn.Range("G10").Value = "c131e10c-ac85-432b-a277-ff971a1a9a8b.d44d1f3c5c326fae8255fd8980e2c32f786f1552.zh-cn.xlf"
$wsZhCn.Range("H10").Value = "2017-02-21 02:59:00"
$wsZhCn.Range("H10").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("I10").Value = "'"
$wsZhCn.Range("I10").Style = "Normal"
$wsZhCn.Range("J10").Value = "'"
$wsZhCn.Range("J10").Style = "Normal"
$wsZhCn.Range("K10").Value = "'"
$wsZhCn.Range("K10").Style = "Normal"
$wsZhCn.Range("L10").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("L10").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("M10").Value = "'"
$wsZhCn.Range("M10").Style = "Normal"
$wsZhCn.Range("N10").Value = "'"
$wsZhCn.Range("N10").Style = "Normal"
$wsZhCn.Range("O10").Value = "True"
$wsZhCn.Range("P10").Value = "'"
$wsZhCn.Range("P10").Style = "Normal"
$wsZhCn.Range("Q10").Value = "False"
$wsZhCn.Range("R10").Value = "'"
$wsZhCn.Range("R10").Style = "Normal"

# ---------------------------------------------------------------------------
# Sheet "de-de": add the matching detail row.
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)
$rowDeDe = $loDeDe.ListRows.Add()

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A10"), "https://github.com/OpenLocalizationTestOrg/ol-test4/blob/c131e10c1e1c131e10c1e1c131e10c1e1c131e10/e2e/c131e10c-ac85-432b-a277-ff971a1a9a8b.md", "", "", "c131e10c-ac85-432b-a277-ff971a1a9a8b.md")
$wsDeDe.Range("B10").Value = ".md"
$wsDeDe.Range("C10").Value = "Ready for handoff"
$wsDeDe.Range("D10").Value = "e2e"
$wsDeDe.Range("E10").Value = "ht"
$wsDeDe.Range("F10").Value = "False"
$wsDeDe.Range("G10").Value = "c131e10c-ac85-432b-a277-ff971a1a9a8b.d44d1f3c5c326fae8255fd8980e2c32f786f1552.de-de.xlf"
$wsDeDe.Range("H10").Value = "2017-02-21 02:59:17"
$wsDeDe.Range("H10").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("I10").Value = "'"
$wsDeDe.Range("I10").Style = "Normal"
$wsDeDe.Range("J10").Value = "'"
$wsDeDe.Range("J10").Style = "Normal"
$wsDeDe.Range("K10").Value = "'"
$wsDeDe.Range("K10").Style = "Normal"
$wsDeDe.Range("L10").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("L10").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("M10").Value = "'"
$wsDeDe.Range("M10").Style = "Normal"
$wsDeDe.Range("N10").Value = "'"
$wsDeDe.Range("N10").Style = "Normal"
$wsDeDe.Range("O10").Value = "True"
$wsDeDe.Range("P10").Value = "'"
$wsDeDe.Range("P10").Style = "Normal"
$wsDeDe.Range("Q10").Value = "False"
$wsDeDe.Range("R10").Value = "'"
$wsDeDe.Range("R10").Style = "Normal"
